# color-text.docx style-sheet update:
#   - switch the East-Asian font used by the doc defaults / Normal /
#     Heading styles from "DejaVu Sans" to "Tahoma"
#   - give the List, Caption and Index styles an explicit complex-script
#     (w:cs) font of "DejaVu Sans" (they previously inherited it)
#
# NOTE: the w:docDefaults / w:rPrDefault block that backs the document's
# overall default run properties is not reachable through the Word
# object model (ActiveDocument.Styles only enumerates the named styles
# actually defined in styles.xml; there is no COM entry point for the
# docDefaults element itself, in this runtime or in real Word
# automation). The style-level edits below are applied via the
# supported Styles/Font surface.

$d = $word.ActiveDocument

# Font.NameFarEast <-> w:rFonts/@w:eastAsia
$d.Styles.Item("Normal").Font.NameFarEast = "Tahoma"
$d.Styles.Item("Heading").Font.NameFarEast = "Tahoma"

# Font.NameBi <-> w:rFonts/@w:cs
$d.Styles.Item("List").Font.NameBi = "DejaVu Sans"
$d.Styles.Item("Caption").Font.NameBi = "DejaVu Sans"
$d.Styles.Item("Index").Font.NameBi = "DejaVu Sans"
